# "Simu MAXSPEED rigidité châssis" - Rectification de la simulation Maxspeed
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MAX SPEED")

# Updated displacement results (X, Y, Z) for the 10 measurement points after
# rectifying the MAX SPEED simulation.
$ws.Range("D6").Value = -0.0867
$ws.Range("E6").Value = -0.00355
$ws.Range("F6").Value = 5.36

$ws.Range("D7").Value = -0.055
$ws.Range("E7").Value = -0.0048
$ws.Range("F7").Value = 5.35

$ws.Range("D8").Value = -0.0564
$ws.Range("E8").Value = -0.000192
$ws.Range("F8").Value = 4.41

$ws.Range("D9").Value = -0.0562
$ws.Range("E9").Value = -0.00169
$ws.Range("F9").Value = 4.41

$ws.Range("D10").Value = -0.0896
$ws.Range("E10").Value = 0.0104
$ws.Range("F10").Value = 3.17

$ws.Range("D11").Value = -0.0589
$ws.Range("E11").Value = 0.00121
$ws.Range("F11").Value = 3.2

$ws.Range("D12").Value = 0.107
$ws.Range("E12").Value = -0.00032
$ws.Range("F12").Value = 1.2

$ws.Range("D13").Value = 0.018
$ws.Range("E13").Value = 0.0056
$ws.Range("F13").Value = 1.2

$ws.Range("D14").Value = -0.016
$ws.Range("E14").Value = 0.0235
$ws.Range("F14").Value = 0.064

$ws.Range("D15").Value = 0.0325
$ws.Range("E15").Value = -0.00108
$ws.Range("F15").Value = 0.576

# The "Notes" annotation cell shifts one column to the right, from I7 to J7.
$ws.Range("I7").Clear()
$ws.Range("J7").Value = "Appliquer les mêmes efforts de la roue droite sur la roue gauche, en inversant la valeur en yrts de la roue droite sur la roue gauche"

# Restore the current selection to reflect where the author left off editing.
[void]$ws.Range("I12").Select()
